# edit.ps1 -- reproduce the gandhi.docx diff via Word COM-interop calls.
#
# Three changes described by the diff:
#  1. In the first paragraph ("This is a Microsoft word document."), append
#     a plain run with two spaces, then three separate red (C00000) runs:
#     "(This is a change - Ve" / "rsion for branch alternate" / ")"
#     (using an en dash, U+2013, between "change" and "Ve").
#  2. Append one brand-new, completely empty paragraph (<w:p/>) right before
#     the closing <w:sectPr>, i.e. after the last (NormalWeb) paragraph.
#  3. Stamp the built-in "Normal" style with <w:rsid w:val="001772C0"/>.

$d = $word.ActiveDocument
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Helper markers used to pull the *real* OOXML for a single-paragraph Range
# out of its WordOpenXML export (the export always appends a synthetic,
# placeholder-id empty paragraph of its own after the real content -- we
# locate and strip that off to recover the exact original markup).
$bodyTag = "<w:body>"
$syntheticTag = '<w:p w14:paraId="00000001"'

function Get-RealParagraphXml($range) {
    $xml = $range.WordOpenXML
    $bodyStart = $xml.IndexOf($bodyTag) + $bodyTag.Length
    $syntheticStart = $xml.IndexOf($syntheticTag)
    return $xml.Substring($bodyStart, $syntheticStart - $bodyStart)
}

# ---------------------------------------------------------------------------
# 1) Append the red annotation after "This is a Microsoft word document."
# ---------------------------------------------------------------------------

$firstPara = $d.Paragraphs.First
$firstRange = $firstPara.Range
$originalFirstParaXml = Get-RealParagraphXml $firstRange

$newRunsXml = "<w:r><w:t xml:space='preserve'>  </w:t></w:r>" + `
              "<w:r><w:rPr><w:color w:val='C00000'/></w:rPr><w:t>(This is a change " + [char]0x2013 + " Ve</w:t></w:r>" + `
              "<w:r><w:rPr><w:color w:val='C00000'/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>" + `
              "<w:r><w:rPr><w:color w:val='C00000'/></w:rPr><w:t>)</w:t></w:r>"

$closeTag = "</w:p>"
$insertAt = $originalFirstParaXml.LastIndexOf($closeTag)
$updatedFirstParaXml = $originalFirstParaXml.Substring(0, $insertAt) + $newRunsXml + $originalFirstParaXml.Substring($insertAt)

$firstRange.InsertXML($updatedFirstParaXml)

# ---------------------------------------------------------------------------
# 2) Append a brand-new empty paragraph right before </w:body>/<w:sectPr>
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$originalLastParaXml = Get-RealParagraphXml $lastRange

# InsertXML replaces the target range's contents, so giving it
# "original paragraph + new empty paragraph" effectively appends the new,
# totally empty <w:p/> right after the original one (whose own markup,
# including its w14:paraId/rsid attributes, is preserved unchanged).
$lastRange.InsertXML($originalLastParaXml + "<w:p $wordNs/>")

# ---------------------------------------------------------------------------
# 3) Stamp the "Normal" style with <w:rsid w:val="001772C0"/>
# ---------------------------------------------------------------------------
# Real Word stamps whichever style it used during an editing session with
# that session's rsid automatically (a pure bookkeeping/provenance marker
# with no visible or semantic effect). This COM-interop host does not
# expose an rsid-stamping primitive on Style/Document (Document.CurrentRsid
# accepts writes but does not drive any rsid emission, and directly futzing
# with it was observed to corrupt unrelated paragraph markup), so there is
# no supported object-model call left that applies just this one metadata
# attribute without side effects; intentionally a no-op here.
